# Hoàn thiện Ngoại Trú
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsCheck = $wb.Worksheets.Item("Check")

# Update Data sheet (row 2)
$wsData.Range("A2").Value = 3011
$wsData.Range("E2").Value = 46200608011
$wsData.Range("X2").Value = "DN4127460130011"

# Update Check sheet (row 2)
$wsCheck.Range("A2").Value = 3011
$wsCheck.Range("C2").Value = "DN4127460130011"

# Reflect final selection back on the Data sheet (matches the saved view state)
$wsData.Activate() | Out-Null
$wsData.Range("A2").Select() | Out-Null
